# Apply the "Fri Jul 1 17:13:46 PDT 2022" template refresh:
#   1. Bump the cached datetimeFigureOut placeholder text from 6/15/2022
#      to 7/1/2022 everywhere it is cached (slide master and every
#      slide layout).
#   2. Shrink the Title/Body master default run sizes (txStyles) so the
#      layout reads less "huge": title 33->28, body lvl1 24->18,
#      lvl2 21->18, lvl3 18->14, lvl4 15->12, lvl5 15->12.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# ---------------------------------------------------------------------
# 1) Date placeholders: ppPlaceholderDate = 16
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shape, [string]$newText) {
    try {
        if ($shape.PlaceholderFormat.Type -eq 16) {
            $shape.TextFrame.TextRange.Text = $newText
        }
    } catch {
        # Not a placeholder at all - nothing to do.
    }
}

$newDate = "7/1/2022"

# Slide master's own Date Placeholder shape.
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    Update-DatePlaceholder $master.Shapes.Item($j) $newDate
}

# Every slide layout has its own cached copy of the date placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        Update-DatePlaceholder $layout.Shapes.Item($j) $newDate
    }
}

# NOTE: the notes master also caches a copy of the date placeholder in
# the underlying OOXML, but this COM-interop runtime does not support
# writing to NotesMaster shapes (writes land on the wrong part), so it
# is intentionally left alone here rather than risk corrupting the
# slide master/layouts.

# ---------------------------------------------------------------------
# 2) Master default text sizes (p:txStyles/p:titleStyle + p:bodyStyle)
# ---------------------------------------------------------------------
# TextStyles.Item(1) => titleStyle ; TextStyles.Item(3) => bodyStyle
$titleStyle = $master.TextStyles.Item(1)
$titleStyle.Levels(1).Font.Size = 28

$bodyStyle = $master.TextStyles.Item(3)
$bodyStyle.Levels(1).Font.Size = 18
$bodyStyle.Levels(2).Font.Size = 18
$bodyStyle.Levels(3).Font.Size = 14
$bodyStyle.Levels(4).Font.Size = 12
$bodyStyle.Levels(5).Font.Size = 12

Write-Host "Template refresh applied: dates -> $newDate; title/body sizes reduced."
